# The commit renames the workbook's single worksheet from
# "11.5 Add Update Work Pak_Part_d" to "Work Pak Part_template"
# (the rest of the diff -- fileVersion/rupBuild, xr:revisionPtr,
# x15ac:absPath, the x16r2 namespace, calcId, etc. -- are just
# application/version fingerprints that Excel itself stamps on save
# and are not user-controllable through the object model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Work Pak Part_template"
